$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns with corrected Diebold-Mariano values
$ws.Range("C2").Value = -0.059581607674993
$ws.Range("D2").Value = 0.9528377021102681

$ws.Range("C3").Value = -0.4437522355803322
$ws.Range("D3").Value = 0.6600330405732944

$ws.Range("C4").Value = -1.642135851615539
$ws.Range("D4").Value = 0.1097806223864255
$ws.Range("G4").Value = "No"

$ws.Range("C5").Value = -0.7807899629386806
$ws.Range("D5").Value = 0.4403286993368472

$ws.Range("C6").Value = -0.6180416264975495
$ws.Range("D6").Value = 0.5406660986549472

$ws.Range("C7").Value = -1.406174053278537
$ws.Range("D7").Value = 0.1687464802114129

$ws.Range("C8").Value = -1.205535294137002
$ws.Range("D8").Value = 0.236319923755633

$ws.Range("C9").Value = -0.9258530667893199
$ws.Range("D9").Value = 0.3610478134327813

$ws.Range("C10").Value = -1.015308123923763
$ws.Range("D10").Value = 0.3171299568483601

$ws.Range("C11").Value = 0.3800366867707264
$ws.Range("D11").Value = 0.7062819384069263
